# "Day 41" activities (MySQL + MCP topics) added to the "Week 6" sheet,
# and the active/selected tab moves from "Week 6" back to "Week 1".

$wb = $excel.ActiveWorkbook

$wk6 = $wb.Worksheets.Item("Week 6")
$wk1 = $wb.Worksheets.Item("Week 1")

# Copy the formatting of the last existing data row (12) down onto the
# new rows (13-21) for columns A (date) and B (day), matching the
# centered/date styles already used in the sheet.
$wk6.Range("A12:B12").Copy() | Out-Null
$wk6.Range("A13:B21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Dates (column A) and "Day 41" label (column B) for every new row
for ($r = 13; $r -le 21; $r++) {
    $wk6.Cells.Item($r, 1).Value = 45941
    $wk6.Cells.Item($r, 2).Value = "Day 41"
}

# Columns C/D are written in the exact order the strings were first
# typed, so new shared-string entries line up with the saved file.
$wk6.Cells.Item(14, 4).Value = "Database AI Agent (using mysql)"
$wk6.Cells.Item(13, 4).Value = "MYSQL installation and setup"
$wk6.Cells.Item(13, 3).Value = "AI Agents"
$wk6.Cells.Item(14, 3).Value = "AI Agents"

$wk6.Cells.Item(15, 3).Value = "MCP"
$wk6.Cells.Item(15, 4).Value = "Model Context Protocol (MCP)"
$wk6.Cells.Item(16, 3).Value = "MCP"
$wk6.Cells.Item(16, 4).Value = "n8n-MCP Agent"
$wk6.Cells.Item(17, 3).Value = "MCP"
$wk6.Cells.Item(17, 4).Value = "MCP Client Server Triggering"
$wk6.Cells.Item(18, 3).Value = "MCP"
$wk6.Cells.Item(18, 4).Value = "Integrating MCPs with IDE (cursor)"
$wk6.Cells.Item(19, 3).Value = "MCP"
$wk6.Cells.Item(19, 4).Value = "Human in the Loop"
$wk6.Cells.Item(20, 3).Value = "MCP"
$wk6.Cells.Item(20, 4).Value = "Fallback AI Model in Agent"
$wk6.Cells.Item(21, 3).Value = "MCP"
$wk6.Cells.Item(21, 4).Value = "Building an n8n flow + MicroSaaS"

# Update the selection on "Week 6" and move the active/selected tab to
# "Week 1".
$wk6.Range("D8").Select() | Out-Null
$wk1.Activate() | Out-Null
$wk1.Range("D46").Select() | Out-Null
